$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: convert B23 from text "3" to a real number 3 (content otherwise unchanged)
$ws.Range("B23").Value = 3

# Row 24: new annotation row, copied/adapted from the old row 23 values
$ws.Range("A24").Value = "Ying Tang"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "3"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "无"
$ws.Range("D24").Value = "DFT"
$ws.Range("E24").Value = "WRI"
$ws.Range("F24").Value = "01473e7f-4f45-41be-bd3f-03c0ff83190e"
$ws.Range("G24").Value = "H1u8fMW0b_annotated.xlsx"
$ws.Range("H24").Value = "The citations are in non-standard format (section 1.2: Kalman (1960))."
